# Creating first version bot
# Adds a new worksheet ("Planilha1") after the existing "Plan1" sheet,
# populates it with a small buy/sell EUR<->USD conversion check, and
# updates the view/selection state on both sheets so the new sheet ends
# up active.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed right after Plan1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planilha1"

# Header row - note the deliberate write order so the workbook's shared
# string table fills up in the same sequence as the source file.
$ws2.Range("A1").Value = "Buy_price"
$ws2.Range("D1").Value = "Sell_price"
$ws2.Range("C1").Value = "qnt_dolar"
$ws2.Range("B1").Value = "Buy euro"
$ws2.Range("E1").Value = "Sell euro"
$ws2.Range("F1").Value = "qnt_dolar"

# Data row
$ws2.Range("A2").Value = 1.104
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Formula = "=B2*A2"
$ws2.Range("D2").Value = 1.11
$ws2.Range("E2").Formula = "=C2/D2"
$ws2.Range("F2").Formula = "=E2*D2"
$ws2.Range("G2").Formula = "=B2-E2"

# Page margins on the new sheet (metric defaults: 1.3 / 2 / 0.8 cm).
$ps = $ws2.PageSetup
$ps.LeftMargin = 36.850393728
$ps.RightMargin = 36.850393728
$ps.TopMargin = 56.692913399999995
$ps.BottomMargin = 56.692913399999995
$ps.HeaderMargin = 22.67716464
$ps.FooterMargin = 22.67716464

# Restore/adjust selections: Plan1 keeps a B1:B4 selection (no longer the
# active tab), Planilha1 becomes the active tab with D2 selected.
$ws1.Range("B1:B4").Select() | Out-Null
$ws2.Range("D2").Select() | Out-Null
